$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Cells.Item($row, 58)
    $cell.Formula = "=""2014-06-11"""
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}
